# Updates to Slides and calculations
# - Electronics: several component quantities zeroed out
# - Dead Mass: quantities/costs tweaked
# - Propulsion: air-density constant in the N3 formula updated; stray helper
#   cell removed
# - Active sheet/selection bookkeeping updated to match the edited session

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Electronics sheet: several component quantities (column B) go to 0
# ---------------------------------------------------------------------------
$wsElec = $wb.Worksheets.Item("Electronics")
$wsElec.Range("B5").Value = 0
$wsElec.Range("B8").Value = 0
$wsElec.Range("B10").Value = 0
$wsElec.Range("B11").Value = 0
$wsElec.Range("B12").Value = 0
$wsElec.Range("B13").Select()

# ---------------------------------------------------------------------------
# Dead Mass sheet: quantity / cost adjustments
# ---------------------------------------------------------------------------
$wsDead = $wb.Worksheets.Item("Dead Mass")
$wsDead.Range("B5").Value = 0
$wsDead.Range("C7").Value = 39
$wsDead.Range("B8").Value = 1.5

# ---------------------------------------------------------------------------
# Propulsion sheet: update the air density constant used in the thrust
# equation, and remove the stray scratch calculation in H14
# ---------------------------------------------------------------------------
$wsProp = $wb.Worksheets.Item("Propulsion")
$wsProp.Range("N3").Formula = "=((9.81*TOTAL!C8/1000)^(3/2))/(4*0.11938*BattVolt*SQRT(2*PI()*0.98))/0.85"
$wsProp.Range("H14").ClearContents()

# ---------------------------------------------------------------------------
# Selection / active-tab bookkeeping
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("TOTAL")
$wsTotal.Range("J8").Select()

$wsProp.Activate()
$wsProp.Range("N4").Select()
